$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 18 and 19) so the sheet shrinks from 19 to 17 rows.
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()

# Row 1: ConceptScheme URI
$ws.Range("B1").Value = "http://purl.org/m4m19/subjects"
$ws.Range("C1").Value = "Main GUPRI (i.e. PID) under which all terms are defined. Preference is to use PURLs or W3IDs as they provide permanent resolvable identifiers."

# Row 2: section header (was a PREFIX row for "gen", now a section title row)
$ws.Range("A2").Value = "Prefixes for controlled vocabularies, schema and ontologies"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""

# Row 3: PREFIX nicest-2-subjects
$ws.Range("A3").Value = "PREFIX"
$ws.Range("B3").Value = "nicest-2-subjects"
$ws.Range("C3").Value = "http://purl.org/m4m19/subjects/"
$ws.Range("D3").Value = "Prefix for our controlled vocabulary since it is rather tedious to write long URLs all the time"

# Row 4: PREFIX skos
$ws.Range("A4").Value = "PREFIX"
$ws.Range("B4").Value = "skos"
$ws.Range("C4").Value = "http://www.w3.org/2004/02/skos/core#"
$ws.Range("D4").Value = "Prefix for SKOS Ontology. This ontology is which our base for defining terms."

# Row 5: PREFIX pav
$ws.Range("A5").Value = "PREFIX"
$ws.Range("B5").Value = "pav"
$ws.Range("C5").Value = "http://purl.org/pav/ "
$ws.Range("D5").Value = "Prefix for Provenance, Authoring and Versioning Ontology which properties such as version and createdOn we will use to describe our controlled vocabulary"

# Row 6: PREFIX dct
$ws.Range("A6").Value = "PREFIX"
$ws.Range("B6").Value = "dct"
$ws.Range("C6").Value = "http://purl.org/dc/terms/ "
$ws.Range("D6").Value = "Prefix for Dublin Core (Terms) Ontology which properties such as title, description, rights, source, etc. we will use to describe our controlled vocabulary as well to define its terms"

# Row 7: section header (was a PREFIX row for "iop", now a section title row)
$ws.Range("A7").Value = "Metadata about vocabulary"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""

# Row 8: dct:title
$ws.Range("A8").Value = "dct:title"
$ws.Range("B8").Value = "NICEST-2 controlled vocabulary of subjects"
$ws.Range("C8").Value = "Title of the vocabulary"

# Row 9: dct:description
$ws.Range("A9").Value = "dct:description"
$ws.Range("B9").Value = "Subjects ..."
$ws.Range("C9").Value = "Description of the controlled vocabulary"

# Row 10: dct:creator
$ws.Range("A10").Value = "dct:creator"
$ws.Range("B10").Value = "https://orcid.org/0000-0002-9381-9693"
$ws.Range("C10").Value = "An ORCID ID of the vocabulary creator, repeat this row as needed."

# Row 11: dct:rights
$ws.Range("A11").Value = "dct:rights"
$ws.Range("B11").Value = "https://spdx.org/licenses/CC0-1.0"
$ws.Range("C11").Value = "License under which the vocabulary is provided"

# Row 12: pav:version
$ws.Range("A12").Value = "pav:version"
$ws.Range("B12").Value = "0.1.0"
$ws.Range("C12").Value = "Vocabulary version"

# Row 13: pav:createdOn
$ws.Range("A13").Value = "pav:createdOn"
$ws.Range("B13").Value = "2021-11-12T12:00+02:00"
$ws.Range("C13").Value = "Date when vocabulary was initially created"

# Row 14: pav:lastUpdatedOn
$ws.Range("A14").Value = "pav:lastUpdatedOn"
$ws.Range("B14").Value = "2021-11-12T12:00+02:00"
$ws.Range("C14").Value = "Date of the last vocabulary update"

# Row 15: section header (was pav:version row, now a section title row)
$ws.Range("A15").Value = "Definition of terms"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""

# Row 16: header row for the term definitions table (unchanged content, was row 18)
$ws.Range("A16").Value = "Identifier"
$ws.Range("B16").Value = "skos:prefLabel@en"
$ws.Range("C16").Value = 'skos:altLabel(separator=",")'
$ws.Range("D16").Value = "skos:definition@en"
$ws.Range("E16").Value = 'dct:source(separator=",")'
$ws.Range("F16").Value = 'skos:broader(separator=",")'
$ws.Range("G16").Value = 'skos:exactMatch(separator=",")'
$ws.Range("H16").Value = 'skos:closeMatch(separator=",")'
$ws.Range("I16").Value = "skos:editorialNote@en"
$ws.Range("J16").Value = 'dct:creator(separator=",")'
$ws.Range("K16").Value = 'dct:contributor(separator=",")'

# Row 17: example term row (was row 19)
$ws.Range("A17").Value = "nicest-2-subjects:TestTerm"
$ws.Range("B17").Value = "TestTerm"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = "Some test term"
